$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # interventions
$ws2 = $wb.Worksheets.Item(2)   # themes
$ws3 = $wb.Worksheets.Item(3)   # metric_effects
$ws4 = $wb.Worksheets.Item(4)   # intervention_effects

# --- Add the new "stages" sheet at the end, and populate it first so that
# its new shared strings (src_intervention_id, dst_intervention_id,
# relation_type) are interned before the other new strings.
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws4)
$ws5.Name = "stages"

$ws5.Range("A1").Value = "src_intervention_id"
$ws5.Range("B1").Value = "dst_intervention_id"
$ws5.Range("C1").Value = "relation_type"

# --- sheet1 (interventions): renumber ids and rename causes
$ws1.Range("A2").Value = 600
$ws1.Range("B2").Value = "one"
$ws1.Range("A3").Value = 601
$ws1.Range("B3").Value = "two"
$ws1.Range("A4").Value = 602
$ws1.Range("B4").Value = "three"
$ws1.Range("A5").Value = 603
$ws1.Range("B5").Value = "four"

# --- sheet2 (themes): add two theme rows (write row 3 before row 2 to match
# the original authoring order captured in the shared-string table)
$ws2.Range("A3").Value = 501
$ws2.Range("B3").Value = "test theme two"
$ws2.Range("A2").Value = 500
$ws2.Range("B2").Value = "test theme one"

# --- sheet3 (metric_effects): add a data row
$ws3.Range("A2").Value = 500
$ws3.Range("B2").Value = "external_wall_area"
$ws3.Range("C2").Value = 600
$ws3.Range("D2").Value = "ratio"
$ws3.Range("E2").Value = 999
$ws3.Range("F2").Value = 9999
$ws3.Range("G2").Value = 1.5

# --- sheet4 (intervention_effects): add a data row
$ws4.Range("A2").Value = 501
$ws4.Range("B2").Value = 601
$ws4.Range("C2").Value = 600
$ws4.Range("D2").Value = "ratio"
$ws4.Range("E2").Value = 0.99
$ws4.Range("F2").Value = 1.99
$ws4.Range("G2").Value = 1.5

# --- sheet5 (stages): add a data row
$ws5.Range("A2").Value = 600
$ws5.Range("B2").Value = 603
$ws5.Range("C2").Value = "prereq"

# --- Column widths matching the updated layout
$ws2.Columns.Item(2).ColumnWidth = 18.71
$ws3.Columns.Item(2).ColumnWidth = 20
$ws3.Columns.Item(3).ColumnWidth = 23.7109375
$ws5.Columns.Item(1).ColumnWidth = 22.29
$ws5.Columns.Item(2).ColumnWidth = 18.71
$ws5.Columns.Item(3).ColumnWidth = 16.57

# --- Selections to mirror the final view state
$ws1.Range("D33").Select()
$ws2.Range("C20").Select()
$ws3.Range("D31").Select()
$ws4.Range("F12").Select()
$ws5.Range("E8").Select()

# --- Make "stages" the active sheet/tab, as in the target workbook
$ws5.Activate()
